# Rebuild the guest list: drop the qr_code/attendance/invitation columns,
# keep guest_name, guest_gender, guest_category, guest_contact, guest_address,
# and replace the single sample row with a full 10-row guest roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so stale columns (F:H) and the old sample row disappear.
$ws.Cells.Clear()

# Header row
$headers = @("guest_name", "guest_gender", "guest_category", "guest_contact", "guest_address")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value2 = $headers[$c]
}

# Guest data rows
$data = @(
    @("Nicholas Arthur",    "Male",   "VIP",     "081234567890", "101st Fake Street"),
    @("Emily Johnson",      "Female", "Regular", "082198765432", "202nd Imaginary Avenue  "),
    @("Michael Smith",      "Male",   "VIP",     "081278945612", "303rd Fantasy Lane  "),
    @("Sophia Brown",       "Female", "VIP",     "082345678901", "404th Fictional Road  "),
    @("James Wilson",       "Male",   "Regular", "081234569876", "505th Mythical Drive  "),
    @("Olivia Martinez",    "Female", "VIP",     "082156734589", "606th Dreamland Blvd  "),
    @("Benjamin Taylor",    "Male",   "Regular", "081298734561", "707th Illusion St  "),
    @("Charlotte Anderson", "Female", "Regular", "082312478956", "808th Fable Court  "),
    @("William Thomas",     "Male",   "VIP",     "081289734502", "909th Storybook Ave  "),
    @("Ava Hernandez",      "Female", "Regular", "082376591234", "1001st Legendary Way  ")
)

$lastRow = 1 + $data.Length

# Column D (guest_contact) must be Text so leading zeros on phone numbers
# survive - set the format on the whole column (header + data) before any
# values are written into it.
$ws.Range("D1:D$lastRow").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value2 = $row[$c]
    }
}

# Column widths to match the refreshed layout.
$ws.Columns.Item(1).ColumnWidth = 20.21875
$ws.Columns.Item(2).ColumnWidth = 15.5546875
$ws.Columns.Item(3).ColumnWidth = 18.109375
$ws.Columns.Item(4).ColumnWidth = 15.88671875
$ws.Columns.Item(5).ColumnWidth = 29.33203125

# The old selection (G5) no longer falls inside the shrunk A1:E11 range -
# move it back to A1, matching the refreshed view.
$ws.Range("A1").Select()

